$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Column layout cleanup: the original first <col> entry spuriously
#    spanned columns A:B (min="1" max="2") even though column B had
#    its own, later, overriding <col> entry. Touching column B's width
#    makes the engine re-split the column-range table so column A's
#    entry narrows to just column A (min="1" max="1"), matching the
#    cleaned-up layout - column A keeps its original 30.7109375 width.
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = $ws.Columns("B").ColumnWidth

# ------------------------------------------------------------------
# 2) Objetivos: row (row 10) gains its Portuguese description text in
#    columns B and C (previously these mistakenly held the teacher's
#    name).
# ------------------------------------------------------------------
$objetivosText = "Proporcionar ao aluno conhecimento básico e compreensão de cinemática e dinâmica do corpo rígido. Desenvolver algumas aplicações práticas com ênfase em problemas bidimensionais. Apresentar conceitos fundamentais e exemplos das vibrações mecânicas."
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# ------------------------------------------------------------------
# 3) Insert a brand-new row 13 so the "Docentes responsáveis:" value
#    ("7797767 - Viktor Pastoukhov") gets its own row instead of being
#    crammed into the "Programa resumido:" row. This pushes every row
#    from the old 13 down to 14 (old 23 -> 24), matching the target
#    layout and row heights automatically.
# ------------------------------------------------------------------
$ws.Rows(13).Insert()

# The inserted row copies column A's bold/label formatting from the
# row above; column A must stay empty on this row, so clear it fully.
$ws.Range("A13").Clear()

# Populate the new row's B/C cells with the teacher info, and give
# them the normal wrap-text (col B) / red wrap-text (col C) look used
# throughout the rest of the sheet instead of the inherited bold style.
$docenteText = "7797767 - Viktor Pastoukhov"
$ws.Range("B13").Value = $docenteText
$ws.Range("C13").Value = $docenteText
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").WrapText = $true
$ws.Range("C13").Font.Color = 255

# ------------------------------------------------------------------
# 4) Programa resumido: row (now row 14) - replace the placeholder
#    "Semestral" text with the real short syllabus summary.
# ------------------------------------------------------------------
$resumoText = "Cinemática do corpo rígidoDinâmica do pontoDinâmica do corpo rígido Introdução às vibrações mecânicas"
$ws.Range("B14").Value = $resumoText
$ws.Range("C14").Value = $resumoText

# ------------------------------------------------------------------
# 5) Programa: row (now row 16) - replace the stray duplicated date
#    text with the real full syllabus content.
# ------------------------------------------------------------------
$programaText = "Cinemática do corpo rígido:Aceleração e velocidade angulares. Vínculo e cinemática do corpo rígido. Rotação em torno de um eixo fixo. Movimento plano e centro de rotação. Composição de movimentos. Composição de movimentos de rotação.Dinâmica do ponto:Princípios da dinâmica do ponto. Teorema da resultante. Teorema da energia cinética para partícula. Teorema da quantidade de movimento.Dinâmica do corpo rígido:Teorema do movimento do baricentro. Teorema da energia cinética para um sistema de partículas. Teorema do momento angular para um sistema de partículas. Teorema da energia cinética para o corpo rígido. Teorema do momento angular para corpo rígido Exercícios de aplicação: problemas bidimensionais. Rotação do corpo rígido, Balanceamento. Movimento de um giroscópio.Introdução às vibrações mecânicas:Vibrações de sistemas mecânicos com um grau de liberdade: livres sem amortecimento, livres com amortecimento, forçadas. Vibrações de sistemas mecânicos com dois e mais graus de liberdade. Exemplos."
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# ------------------------------------------------------------------
# 6) Método: row (now row 19) - replace the stray duplicated teacher
#    name with the real evaluation method text.
# ------------------------------------------------------------------
$metodoText = "A avaliação será composta por duas provas (P1 e P2)."
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# ------------------------------------------------------------------
# 7) Critério: row (now row 20) - the grading-criteria text shifts
#    down from Método's old text to its own row.
# ------------------------------------------------------------------
$criterioText = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# ------------------------------------------------------------------
# 8) Norma de recuperação: row (now row 21) gets the make-up-exam
#    rules text.
# ------------------------------------------------------------------
$normaText = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# ------------------------------------------------------------------
# 9) Bibliografia: row (now row 22) gets the reading list text.
# ------------------------------------------------------------------
$bibliografiaText = "HIBBELER, R.C. Dinâmica - Mecânica para Engenharia. São Paulo: Pearson Brasil, 2011, 12ª ed., 608p. ISBN: 8576058146.BEER, F.P., JOHNSTON Jr., E.R., CLAUSEN, W. E., Mecânica Vetorial para Engenheiros - Dinâmica, 7ª Edição, McGraw-Hill, São Paulo, 2006, 1355 p. FRANÇA, L. N. F., MATSUMURA, A. Z. Mecânica Geral. Edgard Blücher, 2001, 235 p.SOTELO JR., J., FRANÇA, L.N.F., Introdução às vibrações mecânicas, Edgard Blücher, 2006, 168 p. ISBN: 9788521203384.GREENWOOD, D. T. Principles of Dynamics. New York: Prentice-Hall, 2nd ed, 1988, 552 p.TENENBAUM, R. A. Dinâmica. Editora UFRJ, 1997, 756 p.GIACAGLIA, G. E., Mecânica Geral, Editora Campus, Rio de Janeiro, 1982."
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText
